$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find last row of data based on column B (data) having values
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row

# valor_total (C) = quantidade (E) * valor_item (F) for each data row
for ($r = 2; $r -le $lastRow; $r++) {
    $qty = $ws.Cells.Item($r, 5).Value2
    $price = $ws.Cells.Item($r, 6).Value2
    $ws.Cells.Item($r, 3).Value = $qty * $price
}
